$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; this pushes the existing rows 10-30 down to 11-31
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with the new weekly record
$ws.Cells.Item(10, 1).Value = 11
$ws.Cells.Item(10, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(10, 3).Value = "Bíobío"
$ws.Cells.Item(10, 4).Value = 44708
$ws.Cells.Item(10, 5).Value = 8
$ws.Cells.Item(10, 6).Value = 100114007
$ws.Cells.Item(10, 7).Value = "Jengibre"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 50
$ws.Cells.Item(10, 11).Value = 13000
$ws.Cells.Item(10, 12).Value = 14000
$ws.Cells.Item(10, 13).Value = 13600
$ws.Cells.Item(10, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(10, 15).Value = "Perú"
$ws.Cells.Item(10, 16).Value = 1046
$ws.Cells.Item(10, 17).Value = 13
$ws.Cells.Item(10, 18).Value = "Hortaliza"
